$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("References")

# Capture the existing reference URLs (currently in B3 and B4) before
# reshuffling rows to make room for the new descriptive text.
$purdueText = $ws.Range("B3").Value2
$researchgateText = $ws.Range("B4").Value2

# Insert a short description above each reference link, and move the
# second reference (originally B4) down to B6.
$ws.Range("B5").Value = "All pH data from 0 to 14 (Mainly use this one)"
$ws.Range("B2").Value = "Based on Figure 3 pH 5 to 9 (Just an inspirational finding influenced the second KB)"
$ws.Range("B6").Value = $researchgateText
$ws.Range("B4").Value = ""

# Turn the two URL cells into real hyperlinks (keeps their text as the
# displayed text, which is Excel's default when TextToDisplay is omitted).
$ws.Hyperlinks.Add($ws.Range("B3"), "https://www.extension.purdue.edu/extmedia/HO/HO-140-W.pdf")
$ws.Hyperlinks.Add($ws.Range("B6"), "https://www.researchgate.net/figure/The-effect-of-soil-pH-on-nutrient-availability_fig2_277669269")

# Match the recorded selection/view state and make References the active tab.
$ws.Range("G11").Select()
$ws.Activate()

$wb.Save()
